# Applies the "Logging_OKANT" update:
#  - inserts a brand-new case (A 3062-2025) as the new row 2, shifting every
#    existing data row down by one
#  - bumps the "Förändrad" (column C) date to 45701 for every data row
#  - appends three brand-new cases at the bottom of the sheet
#    (A 3002-2025, A 3005-2025, A 2855-2025)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a fresh row at row 2 - pushes the old rows 2..37 down to 3..38
# ---------------------------------------------------------------------------
$ws.Rows(2).Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new row 2 with the new case's data
# ---------------------------------------------------------------------------
$ws.Range("A2").Value2 = "A 3062-2025"

$ws.Range("B2").Value2 = 45678
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"

$ws.Range("C2").Value2 = 45701
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"

$ws.Range("D2").Value2 = "OKÄNT"
$ws.Range("E2").Value2 = "OKÄNT"

$ws.Range("G2").Value2 = 2.4
$ws.Range("H2").Value2 = 3
$ws.Range("I2").Value2 = 2
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 0
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 0
$ws.Range("O2").Value2 = 0
$ws.Range("P2").Value2 = 0
$ws.Range("Q2").Value2 = 3

$ws.Range("R2").Value2 = "Korallrot`nSpindelblomster`nFläcknycklar"
$ws.Range("R2").WrapText = $true

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/artfynd/A 3062-2025 artfynd.xlsx", "A 3062-2025")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/kartor/A 3062-2025 karta.png", "A 3062-2025")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomål/A 3062-2025 FSC-klagomål.docx", "A 3062-2025")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomålsmail/A 3062-2025 FSC-klagomål mail.docx", "A 3062-2025")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsyn/A 3062-2025 tillsynsbegäran.docx", "A 3062-2025")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsynsmail/A 3062-2025 tillsynsbegäran mail.docx", "A 3062-2025")'

# Row 2 keeps the same row height as the rest of the table (the wrapped,
# 3-line species list would otherwise auto-grow the row).
$ws.Rows(2).RowHeight = 15

# ---------------------------------------------------------------------------
# 3) Every data row's "Förändrad" date (column C) moves to 45701
# ---------------------------------------------------------------------------
$ws.Range("C3:C38").Value2 = 45701

# The row that used to be the last row (old row 37, "A 2598-2025") is now
# row 38 - it picks up the standard explicit row height used by every other
# non-final row.
$ws.Rows(38).RowHeight = 15

# ---------------------------------------------------------------------------
# 4) Append three brand-new cases after the existing data
# ---------------------------------------------------------------------------

# --- row 39: A 3002-2025 -----------------------------------------------
$ws.Range("A39").Value2 = "A 3002-2025"
$ws.Range("B39").Value2 = 45678
$ws.Range("B39").NumberFormat = "YYYY-MM-DD"
$ws.Range("C39").Value2 = 45701
$ws.Range("C39").NumberFormat = "YYYY-MM-DD"
$ws.Range("D39").Value2 = "OKÄNT"
$ws.Range("E39").Value2 = "OKÄNT"
$ws.Range("G39").Value2 = 2.4
$ws.Range("H39").Value2 = 0
$ws.Range("I39").Value2 = 0
$ws.Range("J39").Value2 = 0
$ws.Range("K39").Value2 = 0
$ws.Range("L39").Value2 = 0
$ws.Range("M39").Value2 = 0
$ws.Range("N39").Value2 = 0
$ws.Range("O39").Value2 = 0
$ws.Range("P39").Value2 = 0
$ws.Range("Q39").Value2 = 0
$ws.Range("R39").WrapText = $true
$ws.Rows(39).RowHeight = 15

# --- row 40: A 3005-2025 -----------------------------------------------
$ws.Range("A40").Value2 = "A 3005-2025"
$ws.Range("B40").Value2 = 45678
$ws.Range("B40").NumberFormat = "YYYY-MM-DD"
$ws.Range("C40").Value2 = 45701
$ws.Range("C40").NumberFormat = "YYYY-MM-DD"
$ws.Range("D40").Value2 = "OKÄNT"
$ws.Range("E40").Value2 = "OKÄNT"
$ws.Range("G40").Value2 = 2.2
$ws.Range("H40").Value2 = 0
$ws.Range("I40").Value2 = 0
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 0
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = 0
$ws.Range("N40").Value2 = 0
$ws.Range("O40").Value2 = 0
$ws.Range("P40").Value2 = 0
$ws.Range("Q40").Value2 = 0
$ws.Range("R40").WrapText = $true
$ws.Rows(40).RowHeight = 15

# --- row 41: A 2855-2025 (last row - default row height, like the old
#     last row used to have before this update) -------------------------
$ws.Range("A41").Value2 = "A 2855-2025"
$ws.Range("B41").Value2 = 45678
$ws.Range("B41").NumberFormat = "YYYY-MM-DD"
$ws.Range("C41").Value2 = 45701
$ws.Range("C41").NumberFormat = "YYYY-MM-DD"
$ws.Range("D41").Value2 = "OKÄNT"
$ws.Range("E41").Value2 = "OKÄNT"
$ws.Range("G41").Value2 = 1
$ws.Range("H41").Value2 = 0
$ws.Range("I41").Value2 = 0
$ws.Range("J41").Value2 = 0
$ws.Range("K41").Value2 = 0
$ws.Range("L41").Value2 = 0
$ws.Range("M41").Value2 = 0
$ws.Range("N41").Value2 = 0
$ws.Range("O41").Value2 = 0
$ws.Range("P41").Value2 = 0
$ws.Range("Q41").Value2 = 0
$ws.Range("R41").WrapText = $true

Write-Output "done"
